$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44: swap the set ordering (content-only text change)
$ws.Range("E44").Value = "{'str', 'list'}"

# Row 45: Scalpel type flips from 'list' to 'str', and the status flips
# from Neutral to a Loss (red-fill) result.
$ws.Range("E45").Value = "str"
$ws.Range("F45").Value = "Loss"
$ws.Range("F45").Interior.Color = 255

# Row 57: PyType Wins count goes up by one (4 -> 5)
$ws.Range("D57").Value = 5

# Give the new / repurposed summary rows (58-59) the same plain body
# formatting used by every other data row before we fill in values.
$ws.Range("A56:F56").Copy() | Out-Null
$ws.Range("A58:F59").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
# ...except column F, whose source row was flagged "Neutral" (orange);
# the new rows 58/59 aren't status cells, so restore the plain white fill
# used by the rest of the sheet's body rows.
$ws.Range("F58:F59").Interior.Color = 16777215

# Row 58 becomes a new "Scalpel Accuracy:" summary line, and the old
# "Accuracy over PyType" line is pushed down to row 59 with an updated
# value (25 -> 20).
$ws.Range("A58").Value = ""
$ws.Range("B58").Value = ""
$ws.Range("C58").Value = "Scalpel Accuracy:"
$ws.Range("D58").Value = 1000
$ws.Range("E58").Value = ""
$ws.Range("F58").Value = ""

$ws.Range("A59").Value = ""
$ws.Range("B59").Value = ""
$ws.Range("C59").Value = ""
$ws.Range("D59").Value = ""
$ws.Range("E59").Value = "Accuracy over PyType"
$ws.Range("F59").Value = 20
